$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1332.6735
$ws.Range("I15").Value = 1332.6735
$ws.Range("K15").Value = 3998.020500000001
$ws.Range("M15").Value = -3829.020500000001
$ws.Range("H17").Value = 2043.963
$ws.Range("J17").Value = 2037.9615
$ws.Range("L17").Value = 6113.8845
$ws.Range("N17").Value = -6449.8845
$ws.Range("H62").Value = 3712.5
$ws.Range("I62").Value = 3637.5
$ws.Range("J62").Value = 3787.5
$ws.Range("K62").Value = 3637.5
$ws.Range("L62").Value = 3787.5
$ws.Range("M62").Value = -3013.5
$ws.Range("N62").Value = -5035.5
$ws.Range("H65").Value = 3712.5
$ws.Range("I65").Value = 3637.5
$ws.Range("J65").Value = 3787.5
$ws.Range("K65").Value = 18187.5
$ws.Range("L65").Value = 18937.5
$ws.Range("M65").Value = -15067.5
$ws.Range("N65").Value = -25177.5
$ws.Range("H98").Value = 1125.25
$ws.Range("I98").Value = 1100.2667
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1100.2667
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 397.7333000000001
$ws.Range("N98").Value = -4496
$ws.Range("H111").Value = 875.7
$ws.Range("I111").Value = 875.7
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2627.1
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 439.8999999999996
$ws.Range("H122").Value = 1125.25
$ws.Range("I122").Value = 1100.2667
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3300.800099999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -850.8000999999995
$ws.Range("N122").Value = -9400
$ws.Range("H125").Value = 1105.8889
$ws.Range("J125").Value = 1142.8
$ws.Range("L125").Value = 10285.2
$ws.Range("N125").Value = -15205.2
$ws.Range("H138").Value = 2892.7273
$ws.Range("J138").Value = 3039.1785
$ws.Range("L138").Value = 9117.5355
$ws.Range("N138").Value = -19397.5355

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1592.12
$ws.Range("I2").Value = 1447.6316
$ws.Range("J2").Value = 2049.6667
$ws.Range("K2").Value = 1447.6316
$ws.Range("L2").Value = 2049.6667
$ws.Range("M2").Value = -1334.6316
$ws.Range("N2").Value = -2275.6667
$ws.Range("H96").Value = 41663.332
$ws.Range("J96").Value = 41663.332
$ws.Range("L96").Value = 41663.332
$ws.Range("N96").Value = -47155.332
$ws.Range("H116").Value = 1592.12
$ws.Range("I116").Value = 1447.6316
$ws.Range("J116").Value = 2049.6667
$ws.Range("K116").Value = 1447.6316
$ws.Range("L116").Value = 2049.6667
$ws.Range("M116").Value = 846.3684000000001
$ws.Range("N116").Value = -6637.6667
$ws.Range("H132").Value = 3562.2104
$ws.Range("J132").Value = 4109
$ws.Range("L132").Value = 12327
$ws.Range("N132").Value = -17387

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1592.12
$ws.Range("I3").Value = 1447.6316
$ws.Range("J3").Value = 2049.6667
$ws.Range("K3").Value = 1447.6316
$ws.Range("L3").Value = 2049.6667
$ws.Range("M3").Value = -1333.6316
$ws.Range("N3").Value = -2277.6667
$ws.Range("H95").Value = 51249.668
$ws.Range("J95").Value = 51249.668
$ws.Range("L95").Value = 51249.668
$ws.Range("N95").Value = -56741.668
$ws.Range("H105").Value = 3543.35
$ws.Range("I105").Value = 3289.6428
$ws.Range("K105").Value = 3289.6428
$ws.Range("M105").Value = -1542.6428

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.70588
$ws.Range("I7").Value = 113.3
$ws.Range("K7").Value = 113.3
$ws.Range("M7").Value = -0.2999999999999972
$ws.Range("H60").Value = 23333.334
$ws.Range("I60").Value = 12500
$ws.Range("K60").Value = 12500
$ws.Range("M60").Value = -11989
$ws.Range("H74").Value = 61316.332
$ws.Range("J74").Value = 61316.332
$ws.Range("L74").Value = 61316.332
$ws.Range("N74").Value = -63064.332
$ws.Range("H77").Value = 61316.332
$ws.Range("J77").Value = 61316.332
$ws.Range("L77").Value = 183948.996
$ws.Range("N77").Value = -192684.996
$ws.Range("H99").Value = 2904.111
$ws.Range("I99").Value = 3462.7144
$ws.Range("J99").Value = 949
$ws.Range("K99").Value = 3462.7144
$ws.Range("L99").Value = 949
$ws.Range("M99").Value = -1964.7144
$ws.Range("N99").Value = -3945
$ws.Range("H126").Value = 2904.111
$ws.Range("I126").Value = 3462.7144
$ws.Range("J126").Value = 949
$ws.Range("K126").Value = 10388.1432
$ws.Range("L126").Value = 2847
$ws.Range("M126").Value = -7918.143199999999
$ws.Range("N126").Value = -7787
$ws.Range("H133").Value = 109000
$ws.Range("J133").Value = 109000
$ws.Range("L133").Value = 109000
$ws.Range("N133").Value = -114060

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 537.5
$ws.Range("J68").Value = 625
$ws.Range("L68").Value = 1875
$ws.Range("N68").Value = -3497
$ws.Range("H71").Value = 537.5
$ws.Range("J71").Value = 625
$ws.Range("L71").Value = 5625
$ws.Range("N71").Value = -13737
$ws.Range("H107").Value = 651.7143
$ws.Range("I107").Value = 366.14285
$ws.Range("K107").Value = 1098.42855
$ws.Range("M107").Value = 821.5714499999999
$ws.Range("H122").Value = 1107.2727
$ws.Range("I122").Value = 964.2857
$ws.Range("J122").Value = 1357.5
$ws.Range("K122").Value = 8678.5713
$ws.Range("L122").Value = 12217.5
$ws.Range("M122").Value = -6228.5713
$ws.Range("N122").Value = -17117.5
$ws.Range("H131").Value = 48682.2
$ws.Range("J131").Value = 5278.8423
$ws.Range("L131").Value = 15836.5269
$ws.Range("N131").Value = -25916.5269

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 40000
$ws.Range("J15").Value = 40000
$ws.Range("L15").Value = 40000
$ws.Range("N15").Value = -40576
$ws.Range("H34").Value = 55000
$ws.Range("J34").Value = 55000
$ws.Range("L34").Value = 55000
$ws.Range("N34").Value = -55536
$ws.Range("H76").Value = 55000
$ws.Range("J76").Value = 55000
$ws.Range("L76").Value = 55000
$ws.Range("N76").Value = -55630
$ws.Range("H79").Value = 55000
$ws.Range("J79").Value = 55000
$ws.Range("L79").Value = 55000
$ws.Range("N79").Value = -57184
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H126").Value = 3321.6038
$ws.Range("I126").Value = 3011.5925
$ws.Range("K126").Value = 9034.7775
$ws.Range("M126").Value = -6564.7775

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3249.75
$ws.Range("H71").Value = 3249.75
$ws.Range("H122").Value = 3816.6667
$ws.Range("I122").Value = 3264.7144
$ws.Range("K122").Value = 9794.143199999999
$ws.Range("M122").Value = -7344.143199999999
$ws.Range("H132").Value = 2996.9697
$ws.Range("J132").Value = 3208.4443
$ws.Range("L132").Value = 9625.332900000001
$ws.Range("N132").Value = -14685.3329

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9196.083000000001
$ws.Range("I41").Value = 12079.8
$ws.Range("J41").Value = 7136.2856
$ws.Range("K41").Value = 12079.8
$ws.Range("L41").Value = 7136.2856
$ws.Range("M41").Value = -11689.8
$ws.Range("N41").Value = -7916.2856
$ws.Range("H95").Value = 30387.666
$ws.Range("J95").Value = 30387.666
$ws.Range("L95").Value = 30387.666
$ws.Range("N95").Value = -35879.666
$ws.Range("H113").Value = 552.375
$ws.Range("I113").Value = 545.3333
$ws.Range("K113").Value = 1635.9999
$ws.Range("M113").Value = 534.0001
$ws.Range("H126").Value = 3946.923
$ws.Range("I126").Value = 3301
$ws.Range("K126").Value = 9903
$ws.Range("M126").Value = -7433
$ws.Range("H132").Value = 3569
$ws.Range("J132").Value = 4030.8333
$ws.Range("L132").Value = 12092.4999
$ws.Range("N132").Value = -17152.4999
